# Insert a new weekly price record for "Cilantro" (Macroferia Regional de
# Talca) as row 42, pushing the existing rows 42-64 down to 43-65.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 42 (shifts rows 42..64 -> 43..65,
# inheriting formatting/number-format from the row above, same as
# Excel's native "Insert Sheet Rows").
$ws.Rows.Item(42).Insert()

$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Macroferia Regional de Talca"
$ws.Range("C42").Value = "Maule"
$ws.Range("D42").Value = 44830
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = 100112040
$ws.Range("G42").Value = "Cilantro"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 150
$ws.Range("K42").Value = 7000
$ws.Range("L42").Value = 7000
$ws.Range("M42").Value = 7000
$ws.Range("N42").Value = "`$/caja 36 atados"
$ws.Range("O42").Value = "Región del Maule"
$ws.Range("P42").Value = 194
$ws.Range("Q42").Value = 36
$ws.Range("R42").Value = "Hortaliza"
